$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 722, shifting existing rows 722-736 down to 725-739.
$ws.Rows("722:724").Insert()

# Common (unchanged) columns for this product/market block.
$A = 11
$B = "Vega Monumental Concepción"
$C = "Bíobío"
$E = 8
$F = "Fruta"
$G = 100106
$H = "Oleaginosos"
$I = 100106002
$J = "Palta"
$K = "Hass"

# New row 722: Especial
$ws.Cells.Item(722, 1).Value = $A
$ws.Cells.Item(722, 2).Value = $B
$ws.Cells.Item(722, 3).Value = $C
$ws.Cells.Item(722, 4).Value = 44890
$ws.Cells.Item(722, 5).Value = $E
$ws.Cells.Item(722, 6).Value = $F
$ws.Cells.Item(722, 7).Value = $G
$ws.Cells.Item(722, 8).Value = $H
$ws.Cells.Item(722, 9).Value = $I
$ws.Cells.Item(722, 10).Value = $J
$ws.Cells.Item(722, 11).Value = $K
$ws.Cells.Item(722, 12).Value = "Especial"
$ws.Cells.Item(722, 13).Value = 50
$ws.Cells.Item(722, 14).Value = 2900
$ws.Cells.Item(722, 15).Value = 2900
$ws.Cells.Item(722, 16).Value = 2900
$ws.Cells.Item(722, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(722, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(722, 19).Value = 2900
$ws.Cells.Item(722, 20).Value = 1

# New row 723: Primera
$ws.Cells.Item(723, 1).Value = $A
$ws.Cells.Item(723, 2).Value = $B
$ws.Cells.Item(723, 3).Value = $C
$ws.Cells.Item(723, 4).Value = 44890
$ws.Cells.Item(723, 5).Value = $E
$ws.Cells.Item(723, 6).Value = $F
$ws.Cells.Item(723, 7).Value = $G
$ws.Cells.Item(723, 8).Value = $H
$ws.Cells.Item(723, 9).Value = $I
$ws.Cells.Item(723, 10).Value = $J
$ws.Cells.Item(723, 11).Value = $K
$ws.Cells.Item(723, 12).Value = "Primera"
$ws.Cells.Item(723, 13).Value = 100
$ws.Cells.Item(723, 14).Value = 2700
$ws.Cells.Item(723, 15).Value = 2700
$ws.Cells.Item(723, 16).Value = 2700
$ws.Cells.Item(723, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(723, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(723, 19).Value = 2700
$ws.Cells.Item(723, 20).Value = 1

# New row 724: Segunda
$ws.Cells.Item(724, 1).Value = $A
$ws.Cells.Item(724, 2).Value = $B
$ws.Cells.Item(724, 3).Value = $C
$ws.Cells.Item(724, 4).Value = 44890
$ws.Cells.Item(724, 5).Value = $E
$ws.Cells.Item(724, 6).Value = $F
$ws.Cells.Item(724, 7).Value = $G
$ws.Cells.Item(724, 8).Value = $H
$ws.Cells.Item(724, 9).Value = $I
$ws.Cells.Item(724, 10).Value = $J
$ws.Cells.Item(724, 11).Value = $K
$ws.Cells.Item(724, 12).Value = "Segunda"
$ws.Cells.Item(724, 13).Value = 100
$ws.Cells.Item(724, 14).Value = 2500
$ws.Cells.Item(724, 15).Value = 2500
$ws.Cells.Item(724, 16).Value = 2500
$ws.Cells.Item(724, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(724, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(724, 19).Value = 2500
$ws.Cells.Item(724, 20).Value = 1

# Ensure the D column (date) cells use the date style (s="2") like the rest of the column.
$ws.Range("D722:D724").NumberFormat = $ws.Range("D721").NumberFormat

$wb.Save()
